# Update BOC USD rates (auto)
# Adds a new published-value row to "All Published Values" and refreshes
# the "Daily Summary" publishes count to match.

$wb = $excel.ActiveWorkbook

$wsValues = $wb.Worksheets.Item("All Published Values")
$wsSummary = $wb.Worksheets.Item("Daily Summary")

# Append new row 9 on "All Published Values" -- force text storage so
# date-like / numeric-like strings round-trip exactly as text (matches
# the rest of the sheet, which stores everything as text).
$rowRange = $wsValues.Range("A9:J9")
$rowRange.NumberFormat = "@"

$wsValues.Range("A9").Value = "2026-01-02"
$wsValues.Range("B9").Value = "2026-01-02 19:43:40"
$wsValues.Range("C9").Value = "697.85"
$wsValues.Range("D9").Value = "697.85"
$wsValues.Range("E9").Value = "700.79"
$wsValues.Range("F9").Value = "700.79"
$wsValues.Range("G9").Value = "702.88"
$wsValues.Range("H9").Value = "2026/01/02 19:43:40"
$wsValues.Range("I9").Value = "2026-01-02 11:45:15"
$wsValues.Range("J9").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Drop back to the sheet's default (unstyled) look now that the text is
# locked in, so the new row matches the plain styling of rows 2-8.
$rowRange.Style = "Normal"

# Extend the autofilter range to cover the new row
$wsValues.Range("A1:J9").AutoFilter() | Out-Null

# Update the publishes count on "Daily Summary"
$wsSummary.Range("B4").Value = 8
